$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.072.10"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").Value = "3.131.96"
$ws.Range("E3").Value = "  +3.10%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'580.11"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").Value = "'174.77"
$ws.Range("E6").Value = "  +3.72%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.128.69"
$ws.Range("E8").Value = "  +3.06%  "
$ws.Range("D9").Value = "'0.524"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("D11").Value = "'0.155"
$ws.Range("E11").Value = "  +1.49%  "
$ws.Range("D12").Value = "'0.484"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "'37.25"
$ws.Range("E14").Value = "  +1.21%  "
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Value = "3.649.85"
$ws.Range("E16").Value = "  +3.12%  "
$ws.Range("D17").Value = "67.084.76"
$ws.Range("E17").Value = "  +1.08%  "
$ws.Range("D18").Value = "'7.18"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("D19").Value = "3.134.42"
$ws.Range("E19").Value = "  +3.25%  "
$ws.Range("D20").Value = "'16.19"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").Value = "'486.35"
$ws.Range("E21").Value = "  +3.44%  "
$ws.Range("D22").Value = "'0.717"
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("D23").Value = "'7.64"
$ws.Range("E23").Value = "  +3.38%  "
$ws.Range("D24").Value = "'84.30"
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("D25").Value = "'13.20"
$ws.Range("E25").Value = "  +3.11%  "
$ws.Range("E26").Value = "  +2.94%  "
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "'7.98"
$ws.Range("E29").Value = "  -2.99%  "
$ws.Range("E30").Value = "  -2.05%  "
$ws.Range("D31").Value = "'2.68"
$ws.Range("E31").Value = "  +1.42%  "
$ws.Range("D32").Value = "'28.84"
$ws.Range("E32").Value = "  +1.78%  "
$ws.Range("D33").Value = "0.0₃0997"
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").Value = "'5.94"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").Value = "'0.988"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("D38").Value = "'47.45"
$ws.Range("E38").Value = "  -2.25%  "
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("D40").Value = "'50.16"
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("D43").Value = "'8.67"
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("D44").Value = "'2.81"
$ws.Range("E44").Value = "  -1.79%  "
$ws.Range("D45").Value = "2.854.21"
$ws.Range("E45").Value = "  +5.23%  "
$ws.Range("D46").Value = "'384.50"
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").Value = "'136.39"
$ws.Range("E48").Value = "  +1.18%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "'24.97"
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("E51").Value = "  -0.63%  "
